$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 13.441269
$ws.Range("H2").Value = 40.323807
$ws.Range("I2").Value = 0.08973082133481231
$ws.Range("J2").Value = 0.08973082133481232
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.847907
$ws.Range("N2").Value = 32.54372100000001
$ws.Range("O2").Value = 0.1203140543159393
$ws.Range("P2").Value = 0.1203140543159394
$ws.Range("Q2").Value = 145.809636073983
$ws.Range("R2").Value = 1312.286724665847
$ws.Range("S2").Value = 0.01079587891189046
$ws.Range("T2").Value = 0.01079587891189046

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 13.441269
$ws.Range("H3").Value = 40.323807
$ws.Range("I3").Value = 0.08973082133481231
$ws.Range("J3").Value = 0.08973082133481232
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 65.55027766666667
$ws.Range("N3").Value = 196.650833
$ws.Range("O3").Value = 0.7270176327666009
$ws.Range("P3").Value = 0.7270176327666009
$ws.Range("Q3").Value = 881.0789151423591
$ws.Range("R3").Value = 7929.710236281232
$ws.Range("S3").Value = 0.06523588931303805
$ws.Range("T3").Value = 0.06523588931303806

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 13.441269
$ws.Range("H4").Value = 40.323807
$ws.Range("I4").Value = 0.08973082133481231
$ws.Range("J4").Value = 0.08973082133481232
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.4248426666666666
$ws.Range("N4").Value = 1.274528
$ws.Range("O4").Value = 0.004711926795930482
$ws.Range("P4").Value = 0.004711926795930482
$ws.Range("Q4").Value = 5.710424565344
$ws.Range("R4").Value = 51.393821088096
$ws.Range("S4").Value = 0.0004228050614683527
$ws.Range("T4").Value = 0.0004228050614683528

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 13.441269
$ws.Range("H5").Value = 40.323807
$ws.Range("I5").Value = 0.08973082133481231
$ws.Range("J5").Value = 0.08973082133481232
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 13.34022966666667
$ws.Range("N5").Value = 40.020689
$ws.Range("O5").Value = 0.1479563861215291
$ws.Range("P5").Value = 0.1479563861215292
$ws.Range("Q5").Value = 179.309615471447
$ws.Range("R5").Value = 1613.786539243023
$ws.Range("S5").Value = 0.01327624804841543
$ws.Range("T5").Value = 0.01327624804841544

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 54.711535
$ws.Range("H6").Value = 164.134605
$ws.Range("I6").Value = 0.3652416280068742
$ws.Range("J6").Value = 0.3652416280068742
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 10.847907
$ws.Range("N6").Value = 32.54372100000001
$ws.Range("O6").Value = 0.1203140543159393
$ws.Range("P6").Value = 0.1203140543159394
$ws.Range("Q6").Value = 593.5056435072451
$ws.Range("R6").Value = 5341.550791565206
$ws.Range("S6").Value = 0.04394370107046117
$ws.Range("T6").Value = 0.04394370107046119

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 54.711535
$ws.Range("H7").Value = 164.134605
$ws.Range("I7").Value = 0.3652416280068742
$ws.Range("J7").Value = 0.3652416280068742
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 65.55027766666667
$ws.Range("N7").Value = 196.650833
$ws.Range("O7").Value = 0.7270176327666009
$ws.Range("P7").Value = 0.7270176327666009
$ws.Range("Q7").Value = 3586.356310819552
$ws.Range("R7").Value = 32277.20679737597
$ws.Range("S7").Value = 0.2655371037813771
$ws.Range("T7").Value = 0.2655371037813771

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 54.711535
$ws.Range("H8").Value = 164.134605
$ws.Range("I8").Value = 0.3652416280068742
$ws.Range("J8").Value = 0.3652416280068742
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.4248426666666666
$ws.Range("N8").Value = 1.274528
$ws.Range("O8").Value = 0.004711926795930482
$ws.Range("P8").Value = 0.004711926795930482
$ws.Range("Q8").Value = 23.24379442682666
$ws.Range("R8").Value = 209.19414984144
$ws.Range("S8").Value = 0.001720991813994864
$ws.Range("T8").Value = 0.001720991813994864

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 54.711535
$ws.Range("H9").Value = 164.134605
$ws.Range("I9").Value = 0.3652416280068742
$ws.Range("J9").Value = 0.3652416280068742
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.34022966666667
$ws.Range("N9").Value = 40.020689
$ws.Range("O9").Value = 0.1479563861215291
$ws.Range("P9").Value = 0.1479563861215292
$ws.Range("Q9").Value = 729.8644423158717
$ws.Range("R9").Value = 6568.779980842845
$ws.Range("S9").Value = 0.05403983134104098
$ws.Range("T9").Value = 0.054039831341041

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 63.67711
$ws.Range("H10").Value = 191.03133
$ws.Range("I10").Value = 0.4250937452800914
$ws.Range("J10").Value = 0.4250937452800915
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 10.847907
$ws.Range("N10").Value = 32.54372100000001
$ws.Range("O10").Value = 0.1203140543159393
$ws.Range("P10").Value = 0.1203140543159394
$ws.Range("Q10").Value = 690.7633673087701
$ws.Range("R10").Value = 6216.870305778931
$ws.Range("S10").Value = 0.051144751958995
$ws.Range("T10").Value = 0.05114475195899502

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 63.67711
$ws.Range("H11").Value = 191.03133
$ws.Range("I11").Value = 0.4250937452800914
$ws.Range("J11").Value = 0.4250937452800915
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 65.55027766666667
$ws.Range("N11").Value = 196.650833
$ws.Range("O11").Value = 0.7270176327666009
$ws.Range("P11").Value = 0.7270176327666009
$ws.Range("Q11").Value = 4174.052241510877
$ws.Range("R11").Value = 37566.47017359789
$ws.Range("S11").Value = 0.3090506483974205
$ws.Range("T11").Value = 0.3090506483974205

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 63.67711
$ws.Range("H12").Value = 191.03133
$ws.Range("I12").Value = 0.4250937452800914
$ws.Range("J12").Value = 0.4250937452800915
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.4248426666666666
$ws.Range("N12").Value = 1.274528
$ws.Range("O12").Value = 0.004711926795930482
$ws.Range("P12").Value = 0.004711926795930482
$ws.Range("Q12").Value = 27.05275321802667
$ws.Range("R12").Value = 243.47477896224
$ws.Range("S12").Value = 0.00200301060916771
$ws.Range("T12").Value = 0.00200301060916771

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 63.67711
$ws.Range("H13").Value = 191.03133
$ws.Range("I13").Value = 0.4250937452800914
$ws.Range("J13").Value = 0.4250937452800915
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 13.34022966666667
$ws.Range("N13").Value = 40.020689
$ws.Range("O13").Value = 0.1479563861215291
$ws.Range("P13").Value = 0.1479563861215292
$ws.Range("Q13").Value = 849.4672719095967
$ws.Range("R13").Value = 7645.20544718637
$ws.Range("S13").Value = 0.06289533431450817
$ws.Range("T13").Value = 0.06289533431450818

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 17.96553866666667
$ws.Range("H14").Value = 53.896616
$ws.Range("I14").Value = 0.119933805378222
$ws.Range("J14").Value = 0.119933805378222
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 10.847907
$ws.Range("N14").Value = 32.54372100000001
$ws.Range("O14").Value = 0.1203140543159393
$ws.Range("P14").Value = 0.1203140543159394
$ws.Range("Q14").Value = 194.888492660904
$ws.Range("R14").Value = 1753.996433948136
$ws.Range("S14").Value = 0.0144297223745927
$ws.Range("T14").Value = 0.0144297223745927

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 17.96553866666667
$ws.Range("H15").Value = 53.896616
$ws.Range("I15").Value = 0.119933805378222
$ws.Range("J15").Value = 0.119933805378222
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 65.55027766666667
$ws.Range("N15").Value = 196.650833
$ws.Range("O15").Value = 0.7270176327666009
$ws.Range("P15").Value = 0.7270176327666009
$ws.Range("Q15").Value = 1177.646048031237
$ws.Range("R15").Value = 10598.81443228113
$ws.Range("S15").Value = 0.08719399127476518
$ws.Range("T15").Value = 0.0871939912747652

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 17.96553866666667
$ws.Range("H16").Value = 53.896616
$ws.Range("I16").Value = 0.119933805378222
$ws.Range("J16").Value = 0.119933805378222
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.4248426666666666
$ws.Range("N16").Value = 1.274528
$ws.Range("O16").Value = 0.004711926795930482
$ws.Range("P16").Value = 0.004711926795930482
$ws.Range("Q16").Value = 7.632527355249778
$ws.Range("R16").Value = 68.692746197248
$ws.Range("S16").Value = 0.0005651193112995556
$ws.Range("T16").Value = 0.0005651193112995557

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 17.96553866666667
$ws.Range("H17").Value = 53.896616
$ws.Range("I17").Value = 0.119933805378222
$ws.Range("J17").Value = 0.119933805378222
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 13.34022966666667
$ws.Range("N17").Value = 40.020689
$ws.Range("O17").Value = 0.1479563861215291
$ws.Range("P17").Value = 0.1479563861215292
$ws.Range("Q17").Value = 239.6644118987138
$ws.Range("R17").Value = 2156.979707088424
$ws.Range("S17").Value = 0.01774497241756454
$ws.Range("T17").Value = 0.01774497241756455
